$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the "Schematic Callout" part-reference values into column A for rows 3-10
$ws.Range("A4").Value = "D1"
$ws.Range("A6").Value = "C2"
$ws.Range("A7").Value = "C3"
$ws.Range("A8").Value = "R1"
$ws.Range("A9").Value = "JP6"
$ws.Range("A10").Value = "JP1"
$ws.Range("A3").Value = "C4"
$ws.Range("A5").Value = "C1"

# Match the final selection recorded in the saved file
$ws.Range("A5").Select()
